$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 332: one more extra-hospital death recorded (M332: 0 -> 1)
# Columns L/M are formatted as Text ("@"), so a plain numeric .Value write
# gets stored as a text string. Temporarily flip to General, write the
# number, then restore the Text format so the stored style index is
# unchanged (matches the rest of the column).
$ws.Range("M332").NumberFormat = "General"
$ws.Range("M332").Value = 1
$ws.Range("M332").NumberFormat = "@"

# Row 334: new positive cases revised upward (49 -> 70) and one more
# extra-hospital death (M334: 1 -> 3)
$ws.Range("C334").Value = 70

$ws.Range("M334").NumberFormat = "General"
$ws.Range("M334").Value = 3
$ws.Range("M334").NumberFormat = "@"

# Row 335: new positive cases revised upward (22 -> 107)
$ws.Range("C335").Value = 107

# Row 336: this day's figures are now filled in (previously blank)
$ws.Range("C336").Value = 9
$ws.Range("E336").Value = 13
$ws.Range("F336").Value = 9
$ws.Range("G336").Value = 135

$ws.Range("L336").NumberFormat = "General"
$ws.Range("L336").Value = 0
$ws.Range("L336").NumberFormat = "@"

$ws.Range("M336").NumberFormat = "General"
$ws.Range("M336").Value = 0
$ws.Range("M336").NumberFormat = "@"
